# Fix the 2015 column: it had no data but the 2016-2021 values were
# shifted one column to the left (into 2015..2020). Insert the missing
# "no data" marker in column G (2015) for rows 4 and 5 and shift the
# existing values for 2016-2021 (columns G:M) right into H:N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (count of permits) ---
$ws.Range("N4").Value = $ws.Range("M4").Value()
$ws.Range("M4").Value = $ws.Range("L4").Value()
$ws.Range("L4").Value = $ws.Range("K4").Value()
$ws.Range("K4").Value = $ws.Range("J4").Value()
$ws.Range("J4").Value = $ws.Range("I4").Value()
$ws.Range("I4").Value = $ws.Range("H4").Value()
$ws.Range("H4").Value = $ws.Range("G4").Value()
$ws.Range("G4").Value = "…"

# --- Row 5 (total area sq.m) ---
$ws.Range("N5").Value = $ws.Range("M5").Value()
$ws.Range("M5").Value = $ws.Range("L5").Value()
$ws.Range("L5").Value = $ws.Range("K5").Value()
$ws.Range("K5").Value = $ws.Range("J5").Value()
$ws.Range("J5").Value = $ws.Range("I5").Value()
$ws.Range("I5").Value = $ws.Range("H5").Value()
$ws.Range("H5").Value = $ws.Range("G5").Value()
$ws.Range("G5").Value = "…"

# Make sure the new G4:N4 / G5:N5 cells carry the same style as the
# rest of their row (the shift above only moved values, and the
# newly-written cells picked up the default style).
$ws.Range("C4").Copy()
$ws.Range("G4:N4").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("G5:N5").PasteSpecial(-4122)

# --- Clean up now-unused decorative / blank cells ---
$ws.Range("B2:V2").Clear()
$ws.Range("B6").Clear()
$ws.Range("B10").Clear()
